$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text so numeric-looking strings
# like "7.060" or "1.002" are not coerced into numbers / lose trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "23.874.68"
$ws.Range("D3").Value = "1.652.99"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "311.06"
$ws.Range("D7").Value = "0.3892"
$ws.Range("D8").Value = "0.3813"
$ws.Range("D9").Value = "51.59"
$ws.Range("D10").Value = "1.351"
$ws.Range("D11").Value = "1.002"
$ws.Range("D12").Value = "0.08483"
$ws.Range("D13").Value = "24.06"
$ws.Range("D14").Value = "7.060"
$ws.Range("D15").Value = "8.079"
$ws.Range("D16").Value = "0.00001316"
$ws.Range("D17").Value = "1.651.26"
$ws.Range("D18").Value = "94.25"
$ws.Range("D19").Value = "0.07012"
$ws.Range("D20").Value = "19.66"
$ws.Range("D21").Value = "6.980"
$ws.Range("D23").Value = "13.77"
$ws.Range("D24").Value = "23.875.31"
$ws.Range("D25").Value = "2.434"
$ws.Range("D26").Value = "2.971"
$ws.Range("D27").Value = "22.10"
$ws.Range("D28").Value = "154.14"
$ws.Range("D29").Value = "5.439"
$ws.Range("D31").Value = "7.877"
$ws.Range("D32").Value = "2.499"
$ws.Range("D33").Value = "1.843.78"
$ws.Range("D34").Value = "1.017"
$ws.Range("D35").Value = "0.08202"
$ws.Range("D36").Value = "0.02911"
$ws.Range("D37").Value = "6.642"
$ws.Range("D38").Value = "10.82"
$ws.Range("D39").Value = "0.2682"
$ws.Range("D40").Value = "0.09159"
$ws.Range("D43").Value = "1.426"
$ws.Range("D44").Value = "16.44"
$ws.Range("D45").Value = "0.6946"
$ws.Range("D46").Value = "2.455"
$ws.Range("D47").Value = "4.099"
$ws.Range("D49").Value = "0.08305"
$ws.Range("D50").Value = "134.18"
$ws.Range("D51").Value = "1.226"

# Row 41/42: Aptos and TheSandbox swap positions (name, link, price, volume)
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.7590"
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "13.59"

# Restore default (no explicit number format) style on the Price column
# now that the text values are committed, so styling matches the original.
$ws.Range("D2:D51").Style = "Normal"

# --- Volume(1h) (column E) updates - values carry 2 leading/trailing spaces ---
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("E34").Value = "  -4.71%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("E36").Value = "  -5.20%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("E51").Value = "  -3.31%  "
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("E42").Value = "  -1.62%  "
